# "1st changes of mifos to finflux"
# On the "Repayment Schedule" sheet, insert a new (blank) column before
# column N. This pushes the existing "Late" / (blank) / "Outstanding"
# columns one slot to the right (N->O, O->P, P->Q) and widens the used
# range from A1:P15 to A1:Q15.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns("N:N").Insert()

# Make "Repayment Schedule" the active sheet/tab (it was "Transactions"
# before) and move its selection to R5.
$ws.Activate()
$ws.Range("R5").Select() | Out-Null
